# Applies the update described by the commit "Atualizado por script em 05-11-2023 20:45":
#   - Swap home/away (and all odds/dates/url) data between row pairs 18/19, 58/59 and 66/67
#     (the match order on those dates changed but the per-row index/date metadata A:E stays put)
#   - Append two new fixtures as rows 76 and 77 (Basel-Yverdon and Lausanne Ouchy-St. Gallen)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param($rowA, $rowB)

    $dataA = $ws.Range("F${rowA}:V${rowA}").Value()
    $dataB = $ws.Range("F${rowB}:V${rowB}").Value()

    $ws.Range("F${rowA}:V${rowA}").Value = $dataB
    $ws.Range("F${rowB}:V${rowB}").Value = $dataA
}

# --- Swap the mixed-up fixtures back into the right rows ---
Swap-RowData 18 19
Swap-RowData 58 59
Swap-RowData 66 67

# --- Append the two new fixtures, copying row 75's formatting first ---
$ws.Range("A75:V75").Copy()
$ws.Range("A76:V76").PasteSpecial(-4122)
$ws.Range("A75:V75").Copy()
$ws.Range("A77:V77").PasteSpecial(-4122)

function Set-MatchRow {
    param(
        $row, $idx, $home, $homeGoals, $away, $awayGoals,
        $homeOpenOdds, $homeOpenDt, $homeCloseOdds, $homeCloseDt,
        $drawOpenOdds, $drawOpenDt, $drawCloseOdds, $drawCloseDt,
        $awayOpenOdds, $awayOpenDt, $awayCloseOdds, $awayCloseDt,
        $url
    )

    $ws.Range("A$row").Value = $idx
    $ws.Range("B$row").Value = "switzerland"
    $ws.Range("C$row").Value = "super-league"
    $ws.Range("D$row").Value = "2023-2024"
    $ws.Range("E$row").Value = 45235.6875

    $ws.Range("F$row").Value = $home
    $ws.Range("G$row").Value = $homeGoals
    $ws.Range("H$row").Value = $away
    $ws.Range("I$row").Value = $awayGoals

    $ws.Range("J$row").Value = $homeOpenOdds
    $ws.Range("K$row").Value = $homeOpenDt
    $ws.Range("L$row").Value = $homeCloseOdds
    $ws.Range("M$row").Value = $homeCloseDt

    $ws.Range("N$row").Value = $drawOpenOdds
    $ws.Range("O$row").Value = $drawOpenDt
    $ws.Range("P$row").Value = $drawCloseOdds
    $ws.Range("Q$row").Value = $drawCloseDt

    $ws.Range("R$row").Value = $awayOpenOdds
    $ws.Range("S$row").Value = $awayOpenDt
    $ws.Range("T$row").Value = $awayCloseOdds
    $ws.Range("U$row").Value = $awayCloseDt

    $ws.Range("V$row").Value = $url
}

Set-MatchRow 76 75 "Basel" 2 "Yverdon" 1 `
    1.7  "29/10/2023 16:43" 1.83 "05/11/2023 16:05" `
    4.5  "29/10/2023 16:43" 4.06 "05/11/2023 16:20" `
    3.9  "29/10/2023 16:43" 4.11 "05/11/2023 15:44" `
    "https://www.betexplorer.com/football/switzerland/super-league/basel-yverdon/tvYs6YHM/"

Set-MatchRow 77 76 "Lausanne Ouchy" 2 "St. Gallen" 5 `
    4.04 "29/10/2023 16:43" 4.03 "05/11/2023 16:27" `
    4.23 "29/10/2023 16:43" 4.08 "05/11/2023 16:27" `
    1.72 "29/10/2023 16:43" 1.85 "05/11/2023 16:23" `
    "https://www.betexplorer.com/football/switzerland/super-league/lausanne-ouchy-st-gallen/ILzBaWfj/"
